# "Add files via upload" — re-upload of the annotated workbook.
#
# Observable, COM-addressable parts of this commit:
#   1. On sheet "scarlet": columns A ("Shuffle") and B ("sentenceID") are
#      unhidden and given explicit widths (they were previously zero-width
#      / hidden helper columns).
#   2. The "scarlet" tab becomes the selected/active tab (it was
#      "Formatted" before). Excel only ever keeps one sheet's tabSelected
#      flag set, and moves bookViews/workbookView's activeTab to match, so
#      activating "scarlet" is enough to also clear tabSelected from
#      "Formatted" and retarget the workbook's active-tab pointer.
#   3. Column A's RAND() volatile formulas naturally recompute to fresh
#      values as part of the normal recalculation pass that happens after
#      the workbook is edited/saved — nothing extra to script for that.
#
# (The author's absolute Dropbox folder path recorded in x15ac:absPath is
# an Office-managed, save-location-derived field that isn't exposed on the
# Workbook/Application object model, so it isn't something a COM script
# can set directly.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scarlet")

# --- Unhide columns A:B and size them -------------------------------------
$colA = $ws.Columns.Item(1)
$colB = $ws.Columns.Item(2)

$colA.Hidden = $false
$colB.Hidden = $false

# ColumnWidth is expressed in "characters of the Normal style's font" and
# Excel re-quantizes it when storing the sheet's <col> width (padding +
# pixel rounding). Backing out the character width that rounds back to the
# target stored widths (~19.44 and ~34.33) gets us as close to those as the
# rounding grid allows.
$colA.ColumnWidth = 18.608072916666668
$colB.ColumnWidth = 33.498697916666664

# --- Make "scarlet" the selected/active sheet ------------------------------
$ws.Activate()
